$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 from 10 to 15
$ws.Range("B1").Value = 15

# Replace the formula in D1 with a static computed value (600)
$ws.Range("D1").Value = 600

# Update B2 from 105 to 110
$ws.Range("B2").Value = 110

# Update D2 from 5250 to 5500
$ws.Range("D2").Value = 5500
